# Update countries & provincias Spain
# Refreshes case-count data for several countries; some rows are re-sorted
# (the sheet is kept sorted descending by "Casos totales", column B) so
# country names in column A are also rewritten where rows change position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Rows where only the statistics changed (country stays the same)
Set-Row 4   "Estados Unidos"        988928 1768 118783 814684 15143 48 55461
Set-Row 8   "Alemania"              158142 372  114500 37657  2570  9  5985
Set-Row 17  "Paises Bajos"          38245  400  0      33477  905   43 4518
Set-Row 47  "Republica Dominicana"  6293   158  993    5018   144   4  282
Set-Row 105 "Sri Lanka"             581    58   126    448    2     0  7
Set-Row 130 "Maldivas"              226    12   17     209    2     0  0

# Somalia's total rose from 436 to 480, moving it above "Consejo Danes para
# los Refugiados", "Malta" and "Jordania" (stable sort keeps their relative
# order, just pushed down one row each).
Set-Row 109 "Somalia"                              480 44 10 444 2 3 26
Set-Row 110 "Consejo Danes para los Refugiados"     459 17 50 381 0 0 28
Set-Row 111 "Malta"                                 450 2  286 160 1 0 4
Set-Row 112 "Jordania"                              447 0  337 103 5 0 7

# "Belice" and "Granada" are tied on total (18) but Belice now sorts first.
Set-Row 184 "Belice"   18 0 5  11 1 0 2
Set-Row 185 "Granada"  18 0 7  11 4 0 0

# "San Vicente y las Granadinas" total rose from 14 to 15, tying with
# "Santa Lucia"; San Vicente now sorts first.
Set-Row 193 "San Vicente y las Granadinas" 15 1 5  10 0 0 0
Set-Row 194 "Santa Lucia"                  15 0 15 0  0 0 0

# "Islas Turcas y Caicos" total rose from 11 to 12, moving it above
# "Montserrat" and "Burundi".
Set-Row 197 "Islas Turcas y Caicos" 12 1 5 6 0 0 1
Set-Row 198 "Montserrat"            11 0 2 8 1 0 1
Set-Row 199 "Burundi"               11 0 4 6 0 0 1
